# Auto-generated cell update script
# Applies updated market-data values (currentAveragePrice* and LeveProfit* columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Alpha Profits" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (98 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1005003.8
$ws.Range("I11").Value = 1005003.8
$ws.Range("K11").Value = 1005003.8
$ws.Range("M11").Value = -1004863.8
$ws.Range("H17").Value = 2698.889
$ws.Range("J17").Value = 2859.6
$ws.Range("L17").Value = 8578.799999999999
$ws.Range("N17").Value = -8914.799999999999
$ws.Range("H30").Value = 12000
$ws.Range("I30").Value = 12000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 36000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -35899
$ws.Range("N30").ClearContents()
$ws.Range("H41").Value = 2195
$ws.Range("I41").Value = 2195
$ws.Range("K41").Value = 2195
$ws.Range("M41").Value = -1755
$ws.Range("H62").Value = 7003
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 7003
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H86").Value = 4072
$ws.Range("I86").Value = 2612.125
$ws.Range("J86").Value = 5239.9
$ws.Range("K86").Value = 2612.125
$ws.Range("L86").Value = 5239.9
$ws.Range("M86").Value = -1489.125
$ws.Range("N86").Value = -7485.9
$ws.Range("H89").Value = 4072
$ws.Range("I89").Value = 2612.125
$ws.Range("J89").Value = 5239.9
$ws.Range("K89").Value = 13060.625
$ws.Range("L89").Value = 26199.5
$ws.Range("M89").Value = -7444.625
$ws.Range("N89").Value = -37431.5
$ws.Range("H92").Value = 1008.8261
$ws.Range("I92").Value = 1005.4
$ws.Range("K92").Value = 1005.4
$ws.Range("M92").Value = 242.6
$ws.Range("H98").Value = 4602.357
$ws.Range("I98").Value = 4687.154
$ws.Range("J98").Value = 3500
$ws.Range("K98").Value = 4687.154
$ws.Range("L98").Value = 3500
$ws.Range("M98").Value = -3189.154
$ws.Range("N98").Value = -6496
$ws.Range("H103").Value = 597.5714
$ws.Range("J103").Value = 626.125
$ws.Range("L103").Value = 1878.375
$ws.Range("N103").Value = -3050.375
$ws.Range("H106").Value = 2129.4
$ws.Range("I106").Value = 2129.4
$ws.Range("K106").Value = 2129.4
$ws.Range("M106").Value = -1498.4
$ws.Range("H107").Value = 909.7083
$ws.Range("I107").Value = 914.6667
$ws.Range("J107").Value = 894.8333
$ws.Range("K107").Value = 914.6667
$ws.Range("L107").Value = 894.8333
$ws.Range("M107").Value = 1005.3333
$ws.Range("N107").Value = -4734.8333
$ws.Range("H111").Value = 843.3333
$ws.Range("I111").Value = 843.3333
$ws.Range("K111").Value = 2529.9999
$ws.Range("M111").Value = 537.0001000000002
$ws.Range("H113").Value = 4121.3687
$ws.Range("I113").Value = 2347.875
$ws.Range("K113").Value = 2347.875
$ws.Range("M113").Value = 906.125
$ws.Range("H116").Value = 2771.0588
$ws.Range("I116").Value = 2500.6667
$ws.Range("K116").Value = 2500.6667
$ws.Range("M116").Value = 941.3332999999998
$ws.Range("H122").Value = 4602.357
$ws.Range("I122").Value = 4687.154
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 14061.462
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -11611.462
$ws.Range("N122").Value = -15400
$ws.Range("H131").Value = 1831.6666
$ws.Range("I131").Value = 1831.6666
$ws.Range("K131").Value = 5494.9998
$ws.Range("M131").Value = -454.9997999999996
$ws.Range("H137").Value = 1849.878
$ws.Range("I137").Value = 1268.8334
$ws.Range("K137").Value = 3806.5002
$ws.Range("M137").Value = -1256.5002
$ws.Range("H141").Value = 163281.17
$ws.Range("I141").Value = 194737.4
$ws.Range("K141").Value = 584212.2
$ws.Range("M141").Value = -579032.2

# --- Sheet: ARM (47 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2947819.8
$ws.Range("I32").Value = 3339796.5
$ws.Range("J32").Value = 7994.25
$ws.Range("K32").Value = 3339796.5
$ws.Range("L32").Value = 7994.25
$ws.Range("M32").Value = -3339509.5
$ws.Range("N32").Value = -8568.25
$ws.Range("H43").Value = 27569.4
$ws.Range("J43").Value = 28337.75
$ws.Range("L43").Value = 28337.75
$ws.Range("N43").Value = -28963.75
$ws.Range("H61").Value = 1885.2222
$ws.Range("I61").Value = 1563.2858
$ws.Range("K61").Value = 1563.2858
$ws.Range("M61").Value = -1351.2858
$ws.Range("H102").Value = 1195.6666
$ws.Range("I102").Value = 1130.6364
$ws.Range("K102").Value = 1130.6364
$ws.Range("M102").Value = 491.3635999999999
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 54188.5
$ws.Range("J109").Value = 54188.5
$ws.Range("L109").Value = 54188.5
$ws.Range("N109").Value = -56962.5
$ws.Range("H110").Value = 3376.2354
$ws.Range("I110").Value = 1717.7142
$ws.Range("K110").Value = 1717.7142
$ws.Range("M110").Value = 327.2858000000001
$ws.Range("H112").Value = 19478.666
$ws.Range("J112").Value = 19478.666
$ws.Range("L112").Value = 19478.666
$ws.Range("N112").Value = -22432.666
$ws.Range("H132").Value = 2138.3
$ws.Range("I132").Value = 1699.9615
$ws.Range("K132").Value = 5099.8845
$ws.Range("M132").Value = -2569.8845
$ws.Range("H136").Value = 1885.2222
$ws.Range("I136").Value = 1563.2858
$ws.Range("K136").Value = 4689.857400000001
$ws.Range("M136").Value = -2139.857400000001
$ws.Range("H138").Value = 57997
$ws.Range("J138").Value = 57997
$ws.Range("L138").Value = 57997
$ws.Range("N138").Value = -68277

# --- Sheet: BSM (30 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7179.7144
$ws.Range("I86").Value = 3178.0908
$ws.Range("J86").Value = 9769
$ws.Range("K86").Value = 3178.0908
$ws.Range("L86").Value = 9769
$ws.Range("M86").Value = -2055.0908
$ws.Range("N86").Value = -12015
$ws.Range("H89").Value = 7179.7144
$ws.Range("I89").Value = 3178.0908
$ws.Range("J89").Value = 9769
$ws.Range("K89").Value = 15890.454
$ws.Range("L89").Value = 48845
$ws.Range("M89").Value = -10274.454
$ws.Range("N89").Value = -60077
$ws.Range("H94").Value = 20260.47
$ws.Range("I94").Value = 22912.1
$ws.Range("K94").Value = 22912.1
$ws.Range("M94").Value = -22461.1
$ws.Range("H105").Value = 1473.04
$ws.Range("I105").Value = 1326.45
$ws.Range("K105").Value = 1326.45
$ws.Range("M105").Value = 420.55
$ws.Range("H112").Value = 99749
$ws.Range("J112").Value = 99749
$ws.Range("L112").Value = 99749
$ws.Range("N112").Value = -102703
$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360

# --- Sheet: CRP (41 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2311.524
$ws.Range("I58").Value = 1706.9286
$ws.Range("K58").Value = 1706.9286
$ws.Range("M58").Value = -1503.9286
$ws.Range("H99").Value = 2347.2273
$ws.Range("I99").Value = 1908.5
$ws.Range("J99").Value = 3115
$ws.Range("K99").Value = 1908.5
$ws.Range("L99").Value = 3115
$ws.Range("M99").Value = -410.5
$ws.Range("N99").Value = -6111
$ws.Range("H105").Value = 2218.476
$ws.Range("I105").Value = 1524.2142
$ws.Range("J105").Value = 3607
$ws.Range("K105").Value = 1524.2142
$ws.Range("L105").Value = 3607
$ws.Range("M105").Value = 222.7858000000001
$ws.Range("N105").Value = -7101
$ws.Range("H122").Value = 3009.8
$ws.Range("I122").Value = 3279.8
$ws.Range("K122").Value = 9839.400000000001
$ws.Range("M122").Value = -7389.400000000001
$ws.Range("H126").Value = 2347.2273
$ws.Range("I126").Value = 1908.5
$ws.Range("J126").Value = 3115
$ws.Range("K126").Value = 5725.5
$ws.Range("L126").Value = 9345
$ws.Range("M126").Value = -3255.5
$ws.Range("N126").Value = -14285
$ws.Range("H134").Value = 4169496.2
$ws.Range("I134").Value = 1329.0476
$ws.Range("K134").Value = 3987.142800000001
$ws.Range("M134").Value = -1452.142800000001
$ws.Range("H136").Value = 2311.524
$ws.Range("I136").Value = 1706.9286
$ws.Range("K136").Value = 5120.7858
$ws.Range("M136").Value = -2570.7858
$ws.Range("H141").Value = 168014.44
$ws.Range("J141").Value = 182141.5
$ws.Range("L141").Value = 182141.5
$ws.Range("N141").Value = -192501.5

# --- Sheet: CUL (28 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 408
$ws.Range("I22").Value = 99
$ws.Range("J22").Value = 562.5
$ws.Range("K22").Value = 297
$ws.Range("L22").Value = 1687.5
$ws.Range("M22").Value = -128
$ws.Range("N22").Value = -2025.5
$ws.Range("H27").Value = 408
$ws.Range("I27").Value = 99
$ws.Range("J27").Value = 562.5
$ws.Range("K27").Value = 297
$ws.Range("L27").Value = 1687.5
$ws.Range("M27").Value = -195
$ws.Range("N27").Value = -1891.5
$ws.Range("H55").Value = 3472.2727
$ws.Range("I55").Value = 924.25
$ws.Range("J55").Value = 4928.2856
$ws.Range("K55").Value = 2772.75
$ws.Range("L55").Value = 14784.8568
$ws.Range("M55").Value = -2595.75
$ws.Range("N55").Value = -15138.8568
$ws.Range("J102").Value = 4500
$ws.Range("L102").Value = 13500
$ws.Range("N102").Value = -18368
$ws.Range("H104").Value = 3279.353
$ws.Range("J104").Value = 3099
$ws.Range("L104").Value = 9297
$ws.Range("N104").Value = -14539

# --- Sheet: GSM (41 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1120666.4
$ws.Range("J11").Value = 1253749.6
$ws.Range("L11").Value = 1253749.6
$ws.Range("N11").Value = -1254027.6
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H70").Value = 6698.2954
$ws.Range("I70").Value = 4581.5
$ws.Range("K70").Value = 4581.5
$ws.Range("M70").Value = -4311.5
$ws.Range("H73").Value = 6698.2954
$ws.Range("I73").Value = 4581.5
$ws.Range("K73").Value = 4581.5
$ws.Range("M73").Value = -3645.5
$ws.Range("H80").Value = 6741.8887
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 6741.8887
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6741.8887
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -8737.8887
$ws.Range("H83").Value = 6741.8887
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 6741.8887
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 33709.4435
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -43693.4435
$ws.Range("H113").Value = 2746.5
$ws.Range("I113").Value = 2494.6
$ws.Range("J113").Value = 4006
$ws.Range("K113").Value = 2494.6
$ws.Range("L113").Value = 4006
$ws.Range("M113").Value = -324.5999999999999
$ws.Range("N113").Value = -8346
$ws.Range("H132").Value = 1384.3889
$ws.Range("I132").Value = 1259.9412
$ws.Range("K132").Value = 3779.8236
$ws.Range("M132").Value = -1249.8236

# --- Sheet: LTW (30 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4973.409
$ws.Range("I40").Value = 3197.9167
$ws.Range("K40").Value = 3197.9167
$ws.Range("M40").Value = -3061.9167
$ws.Range("H46").Value = 4467.3335
$ws.Range("J46").Value = 5212
$ws.Range("L46").Value = 5212
$ws.Range("N46").Value = -5588
$ws.Range("H61").Value = 2111
$ws.Range("I61").Value = 1619.8334
$ws.Range("J61").Value = 8005
$ws.Range("K61").Value = 1619.8334
$ws.Range("L61").Value = 8005
$ws.Range("M61").Value = -1417.8334
$ws.Range("N61").Value = -8409
$ws.Range("H100").Value = 304539.9
$ws.Range("I100").Value = 379174.5
$ws.Range("K100").Value = 379174.5
$ws.Range("M100").Value = -378633.5
$ws.Range("H110").Value = 28333
$ws.Range("J110").Value = 28333
$ws.Range("L110").Value = 28333
$ws.Range("N110").Value = -36513
$ws.Range("H113").Value = 2111
$ws.Range("I113").Value = 1619.8334
$ws.Range("J113").Value = 8005
$ws.Range("K113").Value = 1619.8334
$ws.Range("L113").Value = 8005
$ws.Range("M113").Value = 550.1666
$ws.Range("N113").Value = -12345

# --- Sheet: WVR (64 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 39999
$ws.Range("J15").Value = 39999
$ws.Range("L15").Value = 39999
$ws.Range("N15").Value = -40575
$ws.Range("H28").Value = 24941.5
$ws.Range("I28").Value = 7017
$ws.Range("K28").Value = 7017
$ws.Range("M28").Value = -6669
$ws.Range("H37").Value = 21508.5
$ws.Range("J37").Value = 19810.2
$ws.Range("L37").Value = 19810.2
$ws.Range("N37").Value = -20216.2
$ws.Range("H52").Value = 28856.857
$ws.Range("I52").Value = 28666.666
$ws.Range("J52").Value = 29998
$ws.Range("K52").Value = 28666.666
$ws.Range("L52").Value = 29998
$ws.Range("M52").Value = -28440.666
$ws.Range("N52").Value = -30450
$ws.Range("H70").Value = 12911.75
$ws.Range("J70").Value = 12928.429
$ws.Range("L70").Value = 12928.429
$ws.Range("N70").Value = -13558.429
$ws.Range("H73").Value = 12911.75
$ws.Range("J73").Value = 12928.429
$ws.Range("L73").Value = 12928.429
$ws.Range("N73").Value = -15112.429
$ws.Range("H74").Value = 14510.2
$ws.Range("J74").Value = 15937.75
$ws.Range("L74").Value = 15937.75
$ws.Range("N74").Value = -17809.75
$ws.Range("H77").Value = 14510.2
$ws.Range("J77").Value = 15937.75
$ws.Range("L77").Value = 47813.25
$ws.Range("N77").Value = -57173.25
$ws.Range("H113").Value = 1639.875
$ws.Range("I113").Value = 840
$ws.Range("K113").Value = 2520
$ws.Range("M113").Value = -350
$ws.Range("H122").Value = 3990.7144
$ws.Range("I122").Value = 3826.0908
$ws.Range("J122").Value = 4594.3335
$ws.Range("K122").Value = 11478.2724
$ws.Range("L122").Value = 13783.0005
$ws.Range("M122").Value = -9028.2724
$ws.Range("N122").Value = -18683.0005
$ws.Range("H126").Value = 3875.111
$ws.Range("I126").Value = 5494.5
$ws.Range("J126").Value = 2579.6
$ws.Range("K126").Value = 16483.5
$ws.Range("L126").Value = 7738.799999999999
$ws.Range("M126").Value = -14013.5
$ws.Range("N126").Value = -12678.8
$ws.Range("H132").Value = 2378.875
$ws.Range("I132").Value = 2559.276
$ws.Range("K132").Value = 7677.828
$ws.Range("M132").Value = -5147.828
$ws.Range("H136").Value = 1978.7084
$ws.Range("I136").Value = 1924.091
$ws.Range("J136").Value = 2579.5
$ws.Range("K136").Value = 5772.272999999999
$ws.Range("L136").Value = 7738.5
$ws.Range("M136").Value = -3222.272999999999
$ws.Range("N136").Value = -12838.5

Write-Host "Applied 372 cell value updates and 7 cell clears."
